$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure newly-created rows (20-26) get the same cell style (s="2": vertical-center + wrap text)
# that all other data rows already use, matching the pre-existing formatting.
for ($r = 20; $r -le 26; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.WrapText = $true
        $cell.VerticalAlignment = -4108
    }
}

$ws.Cells.Item(2,1).Value = 2026
$ws.Cells.Item(2,2).Value = "Informatique"
$ws.Cells.Item(2,3).Value = "L1"
$ws.Cells.Item(2,4).Value = "Info"
$ws.Cells.Item(2,5).Value = "M"
$ws.Cells.Item(2,6).Value = 120

$ws.Cells.Item(3,1).Value = 2026
$ws.Cells.Item(3,2).Value = "Informatique"
$ws.Cells.Item(3,3).Value = "L1"
$ws.Cells.Item(3,4).Value = "Info"
$ws.Cells.Item(3,5).Value = "F"
$ws.Cells.Item(3,6).Value = 80

$ws.Cells.Item(4,1).Value = 2026
$ws.Cells.Item(4,2).Value = "Informatique"
$ws.Cells.Item(4,3).Value = "L1"
$ws.Cells.Item(4,4).Value = "SV"
$ws.Cells.Item(4,5).Value = "M"
$ws.Cells.Item(4,6).Value = 220

$ws.Cells.Item(5,1).Value = 2026
$ws.Cells.Item(5,2).Value = "Informatique"
$ws.Cells.Item(5,3).Value = "L1"
$ws.Cells.Item(5,4).Value = "SV"
$ws.Cells.Item(5,5).Value = "F"
$ws.Cells.Item(5,6).Value = 150

$ws.Cells.Item(6,1).Value = 2026
$ws.Cells.Item(6,2).Value = "Informatique"
$ws.Cells.Item(6,3).Value = "L1"
$ws.Cells.Item(6,4).Value = "Info"
$ws.Cells.Item(6,5).Value = "M"
$ws.Cells.Item(6,6).Value = 120

$ws.Cells.Item(7,1).Value = 2026
$ws.Cells.Item(7,2).Value = "ST"
$ws.Cells.Item(7,3).Value = "L2"
$ws.Cells.Item(7,4).Value = "LIST1"
$ws.Cells.Item(7,5).Value = "M"
$ws.Cells.Item(7,6).Value = 80

$ws.Cells.Item(8,1).Value = 2026
$ws.Cells.Item(8,2).Value = "CH"
$ws.Cells.Item(8,3).Value = "L2"
$ws.Cells.Item(8,4).Value = "LIST2"
$ws.Cells.Item(8,5).Value = "M"
$ws.Cells.Item(8,6).Value = 20

$ws.Cells.Item(9,1).Value = 2026
$ws.Cells.Item(9,2).Value = "SV"
$ws.Cells.Item(9,3).Value = "L2"
$ws.Cells.Item(9,4).Value = "GAST2"
$ws.Cells.Item(9,5).Value = "F"
$ws.Cells.Item(9,6).Value = 14

$ws.Cells.Item(10,1).Value = 2026
$ws.Cells.Item(10,2).Value = "SV"
$ws.Cells.Item(10,3).Value = "L3"
$ws.Cells.Item(10,4).Value = "GST3"
$ws.Cells.Item(10,5).Value = "F"
$ws.Cells.Item(10,6).Value = 12

$ws.Cells.Item(11,1).Value = 2026
$ws.Cells.Item(11,2).Value = "ST"
$ws.Cells.Item(11,3).Value = "L3"
$ws.Cells.Item(11,4).Value = "LIST3"
$ws.Cells.Item(11,5).Value = "F"
$ws.Cells.Item(11,6).Value = 11

$ws.Cells.Item(12,1).Value = 2026
$ws.Cells.Item(12,2).Value = "CH"
$ws.Cells.Item(12,3).Value = "L2"
$ws.Cells.Item(12,4).Value = "LIST2"
$ws.Cells.Item(12,5).Value = "F"
$ws.Cells.Item(12,6).Value = 10

$ws.Cells.Item(13,1).Value = 2026
$ws.Cells.Item(13,2).Value = "Informatique"
$ws.Cells.Item(13,3).Value = "L1"
$ws.Cells.Item(13,4).Value = "LIST1"
$ws.Cells.Item(13,5).Value = "M"
$ws.Cells.Item(13,6).Value = 55

$ws.Cells.Item(14,1).Value = 2026
$ws.Cells.Item(14,2).Value = "Informatique"
$ws.Cells.Item(14,3).Value = "L1"
$ws.Cells.Item(14,4).Value = "Info"
$ws.Cells.Item(14,5).Value = "F"
$ws.Cells.Item(14,6).Value = 80

$ws.Cells.Item(15,1).Value = 2026
$ws.Cells.Item(15,2).Value = "Informatique"
$ws.Cells.Item(15,3).Value = "L1"
$ws.Cells.Item(15,4).Value = "Info"
$ws.Cells.Item(15,5).Value = "M"
$ws.Cells.Item(15,6).Value = 120

$ws.Cells.Item(16,1).Value = 2026
$ws.Cells.Item(16,2).Value = "SV"
$ws.Cells.Item(16,3).Value = "L1"
$ws.Cells.Item(16,4).Value = "Info"
$ws.Cells.Item(16,5).Value = "F"
$ws.Cells.Item(16,6).Value = 80

$ws.Cells.Item(17,1).Value = 2026
$ws.Cells.Item(17,2).Value = "SV"
$ws.Cells.Item(17,3).Value = "L1"
$ws.Cells.Item(17,4).Value = "Info"
$ws.Cells.Item(17,5).Value = "M"
$ws.Cells.Item(17,6).Value = 120

$ws.Cells.Item(18,1).Value = 2026
$ws.Cells.Item(18,2).Value = "Informatique"
$ws.Cells.Item(18,3).Value = "L1"
$ws.Cells.Item(18,4).Value = "Info"
$ws.Cells.Item(18,5).Value = "F"
$ws.Cells.Item(18,6).Value = 80

$ws.Cells.Item(19,1).Value = 2025
$ws.Cells.Item(19,2).Value = "SV"
$ws.Cells.Item(19,3).Value = "L2"
$ws.Cells.Item(19,4).Value = "LIM1"
$ws.Cells.Item(19,5).Value = "M"
$ws.Cells.Item(19,6).Value = 200

$ws.Cells.Item(20,1).Value = 2025
$ws.Cells.Item(20,2).Value = "CH"
$ws.Cells.Item(20,3).Value = "L3"
$ws.Cells.Item(20,4).Value = "LIM2"
$ws.Cells.Item(20,5).Value = "F"
$ws.Cells.Item(20,6).Value = 180

$ws.Cells.Item(21,1).Value = 2025
$ws.Cells.Item(21,2).Value = "CH"
$ws.Cells.Item(21,3).Value = "L2"
$ws.Cells.Item(21,4).Value = "LIM3"
$ws.Cells.Item(21,5).Value = "M"
$ws.Cells.Item(21,6).Value = 150

$ws.Cells.Item(22,1).Value = 2025
$ws.Cells.Item(22,2).Value = "PUYSQ"
$ws.Cells.Item(22,3).Value = "L2"
$ws.Cells.Item(22,4).Value = "LIM4"
$ws.Cells.Item(22,5).Value = "F"
$ws.Cells.Item(22,6).Value = 22

$ws.Cells.Item(23,1).Value = 2025
$ws.Cells.Item(23,2).Value = "CH"
$ws.Cells.Item(23,3).Value = "L3"
$ws.Cells.Item(23,4).Value = "LIM5"
$ws.Cells.Item(23,5).Value = "M"
$ws.Cells.Item(23,6).Value = 32

$ws.Cells.Item(24,1).Value = 2025
$ws.Cells.Item(24,2).Value = "PUYSQ"
$ws.Cells.Item(24,3).Value = "L2"
$ws.Cells.Item(24,4).Value = "LIM6"
$ws.Cells.Item(24,5).Value = "F"
$ws.Cells.Item(24,6).Value = 45

$ws.Cells.Item(25,1).Value = 2025
$ws.Cells.Item(25,2).Value = "CH"
$ws.Cells.Item(25,3).Value = "L2"
$ws.Cells.Item(25,4).Value = "LIM7"
$ws.Cells.Item(25,5).Value = "M"
$ws.Cells.Item(25,6).Value = 45

$ws.Cells.Item(26,1).Value = 2025
$ws.Cells.Item(26,2).Value = "PUYSQ"
$ws.Cells.Item(26,3).Value = "L3"
$ws.Cells.Item(26,4).Value = "LIM8"
$ws.Cells.Item(26,5).Value = "F"
$ws.Cells.Item(26,6).Value = 40

$ws.Range("I10").Select()